$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$wsVentasGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasGrupo.Range("H46").Value = 1762.2
$wsVentasGrupo.Range("I46").Value = 486
$wsVentasGrupo.Range("H54").Value = "2 de 52"
$wsVentasGrupo.Range("I54").Value = "2 de 52"

# --- Sheet "VENTA MENSUAL" ---
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F46").Value = 2419.39
$wsVentaMensual.Range("F54").Value = 58578.21

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$wsCumplimiento.Range("D7").Value = 2742.3
$wsCumplimiento.Range("E7").Value = -1342.3
$wsCumplimiento.Range("F7").Value = 1.958785714285714

$wsCumplimiento.Range("D8").Value = 533.63
$wsCumplimiento.Range("E8").Value = 466.37
$wsCumplimiento.Range("F8").Value = 0.53363

$wsCumplimiento.Range("D19").Value = 60993.14
$wsCumplimiento.Range("E19").Value = 33454.30064517914
$wsCumplimiento.Range("F19").Value = 0.6457892303205917

# Column F width shrinks from 25 to 24 (ColumnWidth ~23.15 rounds to stored width 24)
$wsCumplimiento.Columns("F").ColumnWidth = 23.15
